# Remove the unused "Requests" sheet entirely.
$wb = $excel.ActiveWorkbook
$requests = $wb.Worksheets.Item("Requests")
$requests.Delete() | Out-Null

# --- Employees sheet: add a "Position" column, drop the trailing
# "ActivatedTariff" column, and add two sample employee rows. ---
$emp = $wb.Worksheets.Item("Employees")

# Insert a new column before "Name" (column B) for "Position".
$emp.Columns.Item(2).Insert() | Out-Null
$emp.Range("B1").Value = "Position"

# The old "Balance"/"ActivatedTariff" headers (now shifted to columns
# F:G) are no longer part of the sheet.
$emp.Range("F:G").Delete() | Out-Null

# New data rows.
$emp.Range("A2").Value = 1
$emp.Range("B2").Value = "Director"
$emp.Range("C2:E2").Value = "ACER"

$emp.Range("A3").Value = 2
$emp.Range("B3").Value = "Employee"
$emp.Range("C3:E3").Value = "SAMSUNG"

$emp.Range("E3").Select() | Out-Null

# --- Clients sheet: renumber the Id column and rename the sample
# placeholder values. ---
$cli = $wb.Worksheets.Item("Clients")

$cli.Range("A2").Value = 1
$cli.Range("B2:D2").Value = "ASUS"

$cli.Range("A3").Value = 2
$cli.Range("B3:D3").Value = "BENZ"

$cli.Range("A4").Value = 3
$cli.Range("B4:D4").Value = "BMW"

$cli.Range("F10").Select() | Out-Null
